# Weekly update: insert 2 new price rows for the latest reporting date
# at the top of the data block (row 193), pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 193 (existing rows 193:230 shift to 195:232)
$ws.Rows("193:194").Insert()

# New row 193 - Primera quality
$ws.Cells.Item(193, 1).Value = 9
$ws.Cells.Item(193, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(193, 3).Value = "Metropolitana"
$ws.Cells.Item(193, 4).Value = 44511
$ws.Cells.Item(193, 5).Value = 13
$ws.Cells.Item(193, 6).Value = 100112032
$ws.Cells.Item(193, 7).Value = "Zapallo italiano"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 160
$ws.Cells.Item(193, 11).Value = 9000
$ws.Cells.Item(193, 12).Value = 10000
$ws.Cells.Item(193, 13).Value = 9500
$ws.Cells.Item(193, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(193, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(193, 16).Value = 158
$ws.Cells.Item(193, 17).Value = 60
$ws.Cells.Item(193, 18).Value = "Hortaliza"

# New row 194 - Segunda quality
$ws.Cells.Item(194, 1).Value = 9
$ws.Cells.Item(194, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(194, 3).Value = "Metropolitana"
$ws.Cells.Item(194, 4).Value = 44511
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 6).Value = 100112032
$ws.Cells.Item(194, 7).Value = "Zapallo italiano"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Segunda"
$ws.Cells.Item(194, 10).Value = 61
$ws.Cells.Item(194, 11).Value = 7000
$ws.Cells.Item(194, 12).Value = 7000
$ws.Cells.Item(194, 13).Value = 7000
$ws.Cells.Item(194, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(194, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(194, 16).Value = 70
$ws.Cells.Item(194, 17).Value = 100
$ws.Cells.Item(194, 18).Value = "Hortaliza"
